$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update/insert cell values to their new positions & content ---
$ws.Range('B10').Value = 'Apresentação dos fundamentos da Ciência dos Materiais visando a introdução ao estudo das características microestruturais e das propriedades dos materiais, apresentação e discussão de exemplos práticos, bem como fornecer subsídios para o estudo das demais disciplinas do ciclo profissional.'
$ws.Range('C10').Value = 'Apresentação dos fundamentos da Ciência dos Materiais visando a introdução ao estudo das características microestruturais e das propriedades dos materiais, apresentação e discussão de exemplos práticos, bem como fornecer subsídios para o estudo das demais disciplinas do ciclo profissional.'
$ws.Range('A12').Value = 'Docentes responsáveis:'
$ws.Range('B13').Value = '6495737 - Durval Rodrigues Junior'
$ws.Range('C13').Value = '6495737 - Durval Rodrigues Junior'
$ws.Range('B14').Value = '5983729 - Fernando Vernilli Junior'
$ws.Range('C14').Value = '5983729 - Fernando Vernilli Junior'
$ws.Range('B15').Value = '984972 - Hugo Ricardo Zschommler Sandim'
$ws.Range('C15').Value = '984972 - Hugo Ricardo Zschommler Sandim'
$ws.Range('B16').Value = '7459752 - Maria Ismenia Sodero Toledo Faria'
$ws.Range('C16').Value = '7459752 - Maria Ismenia Sodero Toledo Faria'
$ws.Range('A17').Value = 'Programa resumido:'
$ws.Range('B17').Value = 'Ligação atômica. Estrutura cristalina. Defeitos em cristais e estruturas não-cristalinas. Relação microestrutura-propriedade.'
$ws.Range('C17').Value = 'Ligação atômica. Estrutura cristalina. Defeitos em cristais e estruturas não-cristalinas. Relação microestrutura-propriedade.'
$ws.Range('A18').Value = 'Short syllabus:'
$ws.Range('A19').Value = 'Programa:'
$ws.Range('B19').Value = '1. Estrutura atômica, ligações atômicas: ligação iônica, ligação covalente, ligação metálica, ligação de Van der Waals, interações dipolo-dipolo e pontes de hidrogênio. Ligações atômicas e o coeficiente de expansão linear.2. Estrutura cristalina: os sete sistemas e as quatorze redes de Bravais; estruturas de metais, cerâmicas e polímeros; direções e planos atômicos (notação de Miller), número de coordenação, empacotamento atômico linear a planar, Lei de Bragg e difração de raios-X;3. Defeitos em cristais e em estruturas amorfas: soluções sólidas (intersticiais e substitucionais); defeitos de ponto, defeitos de linha (discordâncias e sua dinâmica: movimentação e interação), defeitos bidimensionais (falhas de empilhamento, contornos de antifase, contornos de alto e de baixo ângulo), sólidos amorfos, vidros metálicos, defeitos tridimensionais (poros, trincas e inclusões).4. Relação microestrutura-propriedade: exemplos práticos e estudos de caso (propriedades mecânicas, elétricas e magnéticas).Em todos os itens, são abordados os aspectos práticos de cada tópico da ementa para ampliar as competências dos alunos, que serão trabalhados com Estudos de Caso.'
$ws.Range('C19').Value = '1. Estrutura atômica, ligações atômicas: ligação iônica, ligação covalente, ligação metálica, ligação de Van der Waals, interações dipolo-dipolo e pontes de hidrogênio. Ligações atômicas e o coeficiente de expansão linear.2. Estrutura cristalina: os sete sistemas e as quatorze redes de Bravais; estruturas de metais, cerâmicas e polímeros; direções e planos atômicos (notação de Miller), número de coordenação, empacotamento atômico linear a planar, Lei de Bragg e difração de raios-X;3. Defeitos em cristais e em estruturas amorfas: soluções sólidas (intersticiais e substitucionais); defeitos de ponto, defeitos de linha (discordâncias e sua dinâmica: movimentação e interação), defeitos bidimensionais (falhas de empilhamento, contornos de antifase, contornos de alto e de baixo ângulo), sólidos amorfos, vidros metálicos, defeitos tridimensionais (poros, trincas e inclusões).4. Relação microestrutura-propriedade: exemplos práticos e estudos de caso (propriedades mecânicas, elétricas e magnéticas).Em todos os itens, são abordados os aspectos práticos de cada tópico da ementa para ampliar as competências dos alunos, que serão trabalhados com Estudos de Caso.'
$ws.Range('A20').Value = 'Syllabus:'
$ws.Range('A21').Value = 'Avaliação:'
$ws.Range('A22').Value = 'Método:'
$ws.Range('B22').Value = 'Esta é uma disciplina de caráter fundamental, exigindo dedicação individual para assimilação das definições e conceitos. Isto envolve leitura concentrada para fixação dos conceitos teóricos e realização de exercícios numéricos. Duas provas escritas (P1 e P2) serão aplicadas e com pesos iguais. O desenvolvimento do aluno ao longo do curso será aferido e estimulado por meio de discussões sobre um dado tema, porém sem a atribuição de nota, por conta da subjetividade envolvida.'
$ws.Range('C22').Value = 'Esta é uma disciplina de caráter fundamental, exigindo dedicação individual para assimilação das definições e conceitos. Isto envolve leitura concentrada para fixação dos conceitos teóricos e realização de exercícios numéricos. Duas provas escritas (P1 e P2) serão aplicadas e com pesos iguais. O desenvolvimento do aluno ao longo do curso será aferido e estimulado por meio de discussões sobre um dado tema, porém sem a atribuição de nota, por conta da subjetividade envolvida.'
$ws.Range('A23').Value = 'Critério:'
$ws.Range('B23').Value = ': A Nota final (NF) será calculada da seguinte maneira: NF = (0,4*P1 +0,4* P2+ 0,2*NT) / 3'
$ws.Range('C23').Value = ': A Nota final (NF) será calculada da seguinte maneira: NF = (0,4*P1 +0,4* P2+ 0,2*NT) / 3'
$ws.Range('A24').Value = 'Norma de recuperação:'
$ws.Range('B24').Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR) / 2'
$ws.Range('C24').Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR) / 2'
$ws.Range('A25').Value = 'Bibliografia:'
$ws.Range('B25').Value = '1. CALLISTER Jr, W.D., RETHWISCH, D.G. Ciência e Engenharia de Materiais: Uma Introdução, 8ª ed., LTC Editora, 2013.2. ASKELAND, D.R., PHULÉ, P.P., Ciência e Engenharia dos Materiais, CENGAGE, São Paulo, 2008.3. SHACKELFORD, J.F., Ciência dos Materiais, 6a. ed., Pearson, 2008.4. PADILHA, A.F., Materiais para Engenharia: Microestrutura e Propriedades, Hemus Editora, 1997.5. PADILHA, A.F., Técnicas de Análise Microestrutural, Ed. Hemus, 1985.6. REED-HILL, R.E., Princípios de Metalurgia Física, Guanabara Dois, 1982.7. BRANDON, D.D., KAPLAN, W.D., Microstructural Characterization of Materials, 1st. ed., Wiley, 1999.8. ASHBY, M.F., JONES, D.R.H., Engenharia de Materiais, Elsevier Editora, 2007.9. ASHBY, M.F., SHERCLIFF, H., CEBON, D., Materials: Engineering, Science, Processing and Design, Butterworth-Heinemann, 2010.'
$ws.Range('C25').Value = '1. CALLISTER Jr, W.D., RETHWISCH, D.G. Ciência e Engenharia de Materiais: Uma Introdução, 8ª ed., LTC Editora, 2013.2. ASKELAND, D.R., PHULÉ, P.P., Ciência e Engenharia dos Materiais, CENGAGE, São Paulo, 2008.3. SHACKELFORD, J.F., Ciência dos Materiais, 6a. ed., Pearson, 2008.4. PADILHA, A.F., Materiais para Engenharia: Microestrutura e Propriedades, Hemus Editora, 1997.5. PADILHA, A.F., Técnicas de Análise Microestrutural, Ed. Hemus, 1985.6. REED-HILL, R.E., Princípios de Metalurgia Física, Guanabara Dois, 1982.7. BRANDON, D.D., KAPLAN, W.D., Microstructural Characterization of Materials, 1st. ed., Wiley, 1999.8. ASHBY, M.F., JONES, D.R.H., Engenharia de Materiais, Elsevier Editora, 2007.9. ASHBY, M.F., SHERCLIFF, H., CEBON, D., Materials: Engineering, Science, Processing and Design, Butterworth-Heinemann, 2010.'
$ws.Range('A26').Value = 'Requisitos:'
$ws.Range('B27').Value = 'LOM3018 -  Introdução à Engenharia de Materiais  (Requisito fraco)
'
$ws.Range('C27').Value = 'LOM3018 -  Introdução à Engenharia de Materiais  (Requisito fraco)
'

# --- Clear stale cells left behind by the row shuffle ---
$ws.Range('B12').ClearContents()
$ws.Range('C12').ClearContents()
$ws.Range('A13').ClearContents()
$ws.Range('A14').ClearContents()
$ws.Range('A15').ClearContents()
$ws.Range('A16').ClearContents()
$ws.Range('B18').ClearContents()
$ws.Range('C18').ClearContents()
$ws.Range('B20').ClearContents()
$ws.Range('C20').ClearContents()

# --- Row heights: reset rows that no longer need a custom height ---
$ws.Rows.Item(12).AutoFit()
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).AutoFit()
$ws.Rows.Item(15).AutoFit()

# --- Row heights: set rows that need a (new/changed) custom height ---
$ws.Rows.Item(19).RowHeight = 120
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 60
$ws.Rows.Item(24).RowHeight = 60
$ws.Rows.Item(25).RowHeight = 120
$ws.Rows.Item(27).RowHeight = 30
